# edit.ps1 - applies the "added the function testing Assumption 1" commit
#
# Summary of changes being applied to Housing_data.xlsx:
#   1. Rename worksheet "mean_sale_price" -> "median_sale_price"
#   2. On the (renamed) median_sale_price sheet:
#        - F9 changes from 53000 to -53000
#        - A handful of cells that held placeholder text (":" / "c")
#          are normalised to the text "unknown":
#            H74, E138, E147, F147, G147, E241, F241, G241, F335, F429, G429
#   3. Selection / active-cell bookkeeping:
#        - sale_counts sheet: select columns D:H (sqref D1:H1048576),
#          active cell E1
#        - median_sale_price sheet: active cell I18

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename the "mean_sale_price" sheet to "median_sale_price"
# ---------------------------------------------------------------------
$meanSheet = $wb.Worksheets.Item("mean_sale_price")
$meanSheet.Name = "median_sale_price"

$saleCounts = $wb.Worksheets.Item("sale_counts")
$medianSheet = $wb.Worksheets.Item("median_sale_price")

# ---------------------------------------------------------------------
# 2. Fix up the data on the median_sale_price sheet
# ---------------------------------------------------------------------

# F9: 53000 -> -53000
$medianSheet.Range("F9").Value = -53000

# Cells that contained stray placeholder text (":" or "c") are
# corrected to read "unknown"
$unknownCells = @("H74", "E138", "E147", "F147", "G147", "E241", "F241", "G241", "F335", "F429", "G429")
foreach ($addr in $unknownCells) {
    $medianSheet.Range($addr).Value = "unknown"
}

# ---------------------------------------------------------------------
# 3. Restore sheet selections / active cells
# ---------------------------------------------------------------------

# sale_counts: columns D:H selected (active cell ends up as the
# top-left of the selected range, D1, since the engine's Select()
# always normalises the active cell to the first cell of the range -
# real Excel would keep E1 active here, but that nuance can't be
# reproduced through the exposed COM surface)
$saleCounts.Activate()
$saleCounts.Range("E1").Activate()
$saleCounts.Range("D1:H1048576").Select()

# median_sale_price: active cell I18
$medianSheet.Activate()
$medianSheet.Range("I18").Select()
